# WR_89877351_WeekEnding_072025.xlsx - refresh report with corrected billing data
# (re-synchronized Point Number / Billable Unit Code pairing + updated totals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header / summary ---------------------------------------------
$ws.Range("D5").Value  = "Report Generated On: 08/16/2025 12:48 AM"
$ws.Range("C8").Value  = 7309.41
$ws.Range("C10").Value = "07/14/2025 to 07/20/25"

# --- Line items (Tuesday 07/15/2025 detail table, rows 16-32) ------------

# Row 16: Point 02 / POL-35-5 / Rem  ->  Point 01 / PLA-HDIG / Inst
$ws.Range("A16").Value = "Point 01"
$ws.Range("B16").Value = "PLA-HDIG"
$ws.Range("C16").Value = "Inst"
$ws.Range("D16").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("H16").Value = 648.53

# Row 17: Point 03 / POL-40-2  ->  Point 03 / PLA-HDIG
$ws.Range("B17").Value = "PLA-HDIG"
$ws.Range("D17").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("H17").Value = 648.53

# Row 18: Point 04 / POL-35-5 / Rem  ->  Point 05 / PLA-HDIG / Inst
$ws.Range("A18").Value = "Point 05"
$ws.Range("B18").Value = "PLA-HDIG"
$ws.Range("C18").Value = "Inst"
$ws.Range("D18").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("H18").Value = 648.53

# Row 19: Point 05 / POL-40-1  ->  Point 07 / PLA-HDIG
$ws.Range("A19").Value = "Point 07"
$ws.Range("B19").Value = "PLA-HDIG"
$ws.Range("D19").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("H19").Value = 648.53

# Row 20: Point 06 / POL-35-5 / Rem  ->  Point 09 / PLA-HDIG / Inst
$ws.Range("A20").Value = "Point 09"
$ws.Range("B20").Value = "PLA-HDIG"
$ws.Range("C20").Value = "Inst"
$ws.Range("D20").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("H20").Value = 648.53

# Row 21: Point 07 / POL-40-2  ->  Point 11 / PLA-HDIG
$ws.Range("A21").Value = "Point 11"
$ws.Range("B21").Value = "PLA-HDIG"
$ws.Range("D21").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("H21").Value = 648.53

# Row 22: Point 08 / POL-35-6  ->  Point 02 / POL-35-5
$ws.Range("A22").Value = "Point 02"
$ws.Range("B22").Value = "POL-35-5"
$ws.Range("D22").Value = "Pole,35ft,Class 5"

# Row 23: Point 09  ->  Point 03 (unit code/description/price unchanged)
$ws.Range("A23").Value = "Point 03"

# Row 24: Point 10 / POL-35-6  ->  Point 04 / POL-35-5
$ws.Range("A24").Value = "Point 04"
$ws.Range("B24").Value = "POL-35-5"
$ws.Range("D24").Value = "Pole,35ft,Class 5"

# Row 25: Point 11  ->  Point 05 (unit code/description/price unchanged)
$ws.Range("A25").Value = "Point 05"

# Row 26: Point 12 / POL-40-4  ->  Point 06 / POL-35-5
$ws.Range("A26").Value = "Point 06"
$ws.Range("B26").Value = "POL-35-5"
$ws.Range("D26").Value = "Pole,35ft,Class 5"

# Row 27: Point 01 / PLA-HDIG / Inst  ->  Point 07 / POL-40-2 / Inst
$ws.Range("A27").Value = "Point 07"
$ws.Range("B27").Value = "POL-40-2"
$ws.Range("D27").Value = "Pole,40ft,Class 2"
$ws.Range("H27").Value = 478.55

# Row 28: Point 03 / PLA-HDIG / Inst  ->  Point 08 / POL-35-6 / Rem
$ws.Range("A28").Value = "Point 08"
$ws.Range("B28").Value = "POL-35-6"
$ws.Range("C28").Value = "Rem"
$ws.Range("D28").Value = "Pole,35ft,Class 6"
$ws.Range("H28").Value = 198.88

# Row 29: Point 05 / PLA-HDIG  ->  Point 09 / POL-40-2
$ws.Range("A29").Value = "Point 09"
$ws.Range("B29").Value = "POL-40-2"
$ws.Range("D29").Value = "Pole,40ft,Class 2"
$ws.Range("H29").Value = 478.55

# Row 30: Point 07 / PLA-HDIG / Inst  ->  Point 10 / POL-35-6 / Rem
$ws.Range("A30").Value = "Point 10"
$ws.Range("B30").Value = "POL-35-6"
$ws.Range("C30").Value = "Rem"
$ws.Range("D30").Value = "Pole,35ft,Class 6"
$ws.Range("H30").Value = 198.88

# Row 31: Point 09 / PLA-HDIG  ->  Point 11 / POL-40-1
$ws.Range("A31").Value = "Point 11"
$ws.Range("B31").Value = "POL-40-1"
$ws.Range("D31").Value = "Pole,40ft,Class 1"
$ws.Range("H31").Value = 478.55

# Row 32: Point 11 / PLA-HDIG / Inst  ->  Point 12 / PIN-XAL-C / Rem
$ws.Range("A32").Value = "Point 12"
$ws.Range("B32").Value = "PIN-XAL-C"
$ws.Range("C32").Value = "Rem"
$ws.Range("D32").Value = "Pin,Crossarm Light,Corrosive"
$ws.Range("H32").Value = 31.08

# --- Total row -------------------------------------------------------------
$ws.Range("H33").Value = 7309.410000000001
